$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the 4 new report columns (VERBA, STATUS, FORNECEDOR, CATEGORIA) ---
# New header labels in row 2, columns R:U (18-21), continuing after the
# existing "PERIODO ENTREGA" header in Q2.
$ws.Range("R2").Value = "VERBA"
$ws.Range("S2").Value = "STATUS"
$ws.Range("T2").Value = "FORNECEDOR"
$ws.Range("U2").Value = "CATEGORIA"

# Match the header formatting (navy fill / white bold text) used by the
# rest of row 2 by copying the format from the last existing header cell.
$ws.Range("Q2").Copy()
$ws.Range("R2:U2").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Column widths for the 4 new columns.
$ws.Columns.Item(18).ColumnWidth = 19 + 2/3
$ws.Columns.Item(19).ColumnWidth = 28 + 5/6
$ws.Columns.Item(20).ColumnWidth = 27 + 2/3
$ws.Columns.Item(21).ColumnWidth = 32

# Update the view to match the edited state: scrolled so column J is the
# left-most visible column, with R3 the active/selected cell.
$ws.Range("J1").Select()
$excel.ActiveWindow.ScrollColumn = 10
$ws.Range("R3").Select()
